$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from an existing date cell (G3) onto the
# two new date cells before writing values into them.
$ws.Range("G3").Copy()
$ws.Range("G7:G8").PasteSpecial(-4122)

# Row 7
$ws.Cells.Item(7, 1).Value = 9850.93
$ws.Cells.Item(7, 2).Value = 10217.75
$ws.Cells.Item(7, 3).Value = 296.89
$ws.Cells.Item(7, 4).Value = 286.23
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = -3.59
$ws.Cells.Item(7, 7).Value = 42607.884340277778
$ws.Cells.Item(7, 8).Value = $false

# Row 8
$ws.Cells.Item(8, 1).Value = 9974.07
$ws.Cells.Item(8, 2).Value = 9850.93
$ws.Cells.Item(8, 3).Value = 286.39
$ws.Cells.Item(8, 4).Value = 282.82
$ws.Cells.Item(8, 5).Value = $true
$ws.Cells.Item(8, 6).Value = -1.25
$ws.Cells.Item(8, 7).Value = 42608.616412037038
$ws.Cells.Item(8, 8).Value = $true
